# Bump the published version of the term ValueSet and refresh its date,
# matching commit "Added 1.1.0 of term".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: "Version" | "1.0.0" -> "1.1.0"
$ws.Range("B3").Value = "1.1.0"

# Row 8: "Date" | "2023-06-07T11:52:14+02:00" -> "2023-07-10T23:08:03+02:00"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
